$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02743666666666666
$ws.Range("H2").Value = 0.08231
$ws.Range("I2").Value = 0.007366285056527356
$ws.Range("J2").Value = 0.007366285056527356
$ws.Range("M2").Value = 1.302860333333333
$ws.Range("N2").Value = 3.908581
$ws.Range("O2").Value = 0.9669439908960468
$ws.Range("P2").Value = 0.9669439908960467
$ws.Range("Q2").Value = 0.03574614467888888
$ws.Range("R2").Value = 0.32171530211
$ws.Range("S2").Value = 0.007122785070636473
$ws.Range("T2").Value = 0.007122785070636472

$ws.Range("G3").Value = 0.02743666666666666
$ws.Range("H3").Value = 0.08231
$ws.Range("I3").Value = 0.007366285056527356
$ws.Range("J3").Value = 0.007366285056527356
$ws.Range("O3").Value = 0.008324674682103805
$ws.Range("P3").Value = 0.008324674682103805
$ws.Range("Q3").Value = 0.0003077479444444444
$ws.Range("R3").Value = 0.0027697315
$ws.Range("S3").Value = 0.00006132192671123288
$ws.Range("T3").Value = 0.00006132192671123288

$ws.Range("G4").Value = 0.02743666666666666
$ws.Range("H4").Value = 0.08231
$ws.Range("I4").Value = 0.007366285056527356
$ws.Range("J4").Value = 0.007366285056527356
$ws.Range("M4").Value = 0.033323
$ws.Range("N4").Value = 0.099969
$ws.Range("O4").Value = 0.02473133442184949
$ws.Range("P4").Value = 0.02473133442184949
$ws.Range("Q4").Value = 0.0009142720433333333
$ws.Range("R4").Value = 0.00822844839
$ws.Range("S4").Value = 0.0001821780591796505
$ws.Range("T4").Value = 0.0001821780591796505

$ws.Range("G5").Value = 3.368329
$ws.Range("I5").Value = 0.9043398704228307
$ws.Range("J5").Value = 0.9043398704228307
$ws.Range("M5").Value = 1.302860333333333
$ws.Range("N5").Value = 3.908581
$ws.Range("O5").Value = 0.9669439908960468
$ws.Range("P5").Value = 0.9669439908960467
$ws.Range("Q5").Value = 4.388462243716333
$ws.Range("R5").Value = 39.496160193447
$ws.Range("S5").Value = 0.8744460034330658
$ws.Range("T5").Value = 0.8744460034330657

$ws.Range("G6").Value = 3.368329
$ws.Range("I6").Value = 0.9043398704228307
$ws.Range("J6").Value = 0.9043398704228307
$ws.Range("O6").Value = 0.008324674682103805
$ws.Range("P6").Value = 0.008324674682103805
$ws.Range("Q6").Value = 0.03778142361666666
$ws.Range("R6").Value = 0.34003281255
$ws.Range("S6").Value = 0.007528335223325974
$ws.Range("T6").Value = 0.007528335223325974

$ws.Range("G7").Value = 3.368329
$ws.Range("I7").Value = 0.9043398704228307
$ws.Range("J7").Value = 0.9043398704228307
$ws.Range("M7").Value = 0.033323
$ws.Range("N7").Value = 0.099969
$ws.Range("O7").Value = 0.02473133442184949
$ws.Range("P7").Value = 0.02473133442184949
$ws.Range("Q7").Value = 0.112242827267
$ws.Range("R7").Value = 1.010185445403
$ws.Range("S7").Value = 0.02236553176643906
$ws.Range("T7").Value = 0.02236553176643906

$ws.Range("G8").Value = 0.3288616666666667
$ws.Range("H8").Value = 0.9865849999999999
$ws.Range("I8").Value = 0.08829384452064198
$ws.Range("J8").Value = 0.08829384452064198
$ws.Range("M8").Value = 1.302860333333333
$ws.Range("N8").Value = 3.908581
$ws.Range("O8").Value = 0.9669439908960468
$ws.Range("P8").Value = 0.9669439908960467
$ws.Range("Q8").Value = 0.4284608206538889
$ws.Range("R8").Value = 3.856147385885
$ws.Range("S8").Value = 0.08537520239234461
$ws.Range("T8").Value = 0.0853752023923446

$ws.Range("G9").Value = 0.3288616666666667
$ws.Range("H9").Value = 0.9865849999999999
$ws.Range("I9").Value = 0.08829384452064198
$ws.Range("J9").Value = 0.08829384452064198
$ws.Range("O9").Value = 0.008324674682103805
$ws.Range("P9").Value = 0.008324674682103805
$ws.Range("Q9").Value = 0.003688731694444444
$ws.Range("R9").Value = 0.03319858525
$ws.Range("S9").Value = 0.0007350175320665981
$ws.Range("T9").Value = 0.0007350175320665981

$ws.Range("G10").Value = 0.3288616666666667
$ws.Range("H10").Value = 0.9865849999999999
$ws.Range("I10").Value = 0.08829384452064198
$ws.Range("J10").Value = 0.08829384452064198
$ws.Range("M10").Value = 0.033323
$ws.Range("N10").Value = 0.099969
$ws.Range("O10").Value = 0.02473133442184949
$ws.Range("P10").Value = 0.02473133442184949
$ws.Range("Q10").Value = 0.01095865731833333
$ws.Range("R10").Value = 0.098627915865
$ws.Range("S10").Value = 0.00218362459623078
$ws.Range("T10").Value = 0.00218362459623078

